$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Daily update: append the next day's row of data (row 36)
$ws.Range("A36").Value = 45985
$ws.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B36").Value = 83
$ws.Range("C36").Value = 89
$ws.Range("D36").Value = 88
